$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price-column updates are plain numeric strings (e.g. "1.00", "0.489").
# Force Text format *before* writing the value so Excel keeps them as literal
# text (matching the source inlineStr cells) instead of silently coercing them
# into numbers (which would also strip meaningful trailing zeros).
$textFormatCells = "D5", "D6", "D9", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D25", "D28", "D29", "D30", "D31", "D32", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D47", "D50", "D51"
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed coin values scraped for this run.
$ws.Range("D2").Value = '64.529.28'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '3.507.04'
$ws.Range("E3").Value = '  -1.85%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '585.89'
$ws.Range("E5").Value = '  -2.34%  '
$ws.Range("D6").Value = '132.53'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("D7").Value = '3.506.99'
$ws.Range("E7").Value = '  -1.81%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").Value = '0.389'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").Value = '4.103.77'
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").Value = '27.87'
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '3.507.77'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").Value = '64.526.63'
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = '9.99'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '14.26'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '5.70'
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = '389.57'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").Value = '0.578'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '3.650.05'
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("D25").Value = '74.11'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -3.68%  '
$ws.Range("D28").Value = '1.55'
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.42'
$ws.Range("E30").Value = '  -6.17%  '
$ws.Range("D31").Value = '2.26'
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -5.40%  '
$ws.Range("D33").Value = '3.514.47'
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '23.99'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '5.25'
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.59'
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").Value = '6.95'
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").Value = '0.0811'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '0.813'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").Value = '26.14'
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '42.19'
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("D47").Value = '4.40'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '2.466.61'
$ws.Range("E49").Value = '  +1.06%  '
$ws.Range("D50").Value = '6.90'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("D51").Value = '0.897'
$ws.Range("E51").Value = '  +3.12%  '
